# The deck ships two theme parts (ppt/theme/theme1.xml and
# ppt/theme/theme2.xml). theme2.xml is the one actually wired up as the
# presentation's live theme (slide master + ActivePresentation), so it is
# the only theme reachable through the PowerPoint object model. The edit
# swaps the "Office Theme" colour palette that used to live in theme1.xml
# into the live theme, replacing the "Integral" / "Red Violet" palette it
# had before.
#
# ThemeColorScheme.Colors/Item indices follow the fixed clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB is supplied as a packed COM colour (R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0        # dk1      000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink 954F72
